$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 385.8125
$ws.Range("I2").Value = 193.5
$ws.Range("J2").Value = 578.125
$ws.Range("K2").Value = 193.5
$ws.Range("L2").Value = 578.125
$ws.Range("M2").Value = -80.5
$ws.Range("N2").Value = -804.125
$ws.Range("H17").Value = 983.4375
$ws.Range("J17").Value = 983.4375
$ws.Range("L17").Value = 2950.3125
$ws.Range("N17").Value = -3286.3125
$ws.Range("H40").Value = 2471.1428
$ws.Range("J40").Value = 1913.7142
$ws.Range("L40").Value = 1913.7142
$ws.Range("N40").Value = -2263.7142
$ws.Range("H63").Value = 45271
$ws.Range("J63").Value = 45271
$ws.Range("L63").Value = 45271
$ws.Range("N63").Value = -46519
$ws.Range("H66").Value = 45271
$ws.Range("J66").Value = 45271
$ws.Range("L66").Value = 135813
$ws.Range("N66").Value = -142053
$ws.Range("H129").Value = 872.6727
$ws.Range("J129").Value = 1061.439
$ws.Range("L129").Value = 3184.317
$ws.Range("N129").Value = -13184.317
$ws.Range("H132").Value = 5655263.5
$ws.Range("I132").Value = 7755684.5
$ws.Range("K132").Value = 23267053.5
$ws.Range("M132").Value = -23264523.5
$ws.Range("H138").Value = 558691.4
$ws.Range("I138").Value = 869.9722
$ws.Range("J138").Value = 1036823.94
$ws.Range("K138").Value = 2609.9166
$ws.Range("L138").Value = 3110471.82
$ws.Range("M138").Value = 2530.0834
$ws.Range("N138").Value = -3120751.82

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4326.955
$ws.Range("I32").Value = 3931.0847
$ws.Range("K32").Value = 3931.0847
$ws.Range("M32").Value = -3644.0847
$ws.Range("H45").Value = 1373.5555
$ws.Range("I45").Value = 1346.1428
$ws.Range("J45").Value = 1469.5
$ws.Range("K45").Value = 1346.1428
$ws.Range("L45").Value = 1469.5
$ws.Range("M45").Value = -969.1428000000001
$ws.Range("N45").Value = -2223.5
$ws.Range("H61").Value = 20000942
$ws.Range("I61").Value = 22223018
$ws.Range("J61").Value = 2262.8
$ws.Range("K61").Value = 22223018
$ws.Range("L61").Value = 2262.8
$ws.Range("M61").Value = -22222806
$ws.Range("N61").Value = -2686.8
$ws.Range("H74").Value = 1017.9286
$ws.Range("I74").Value = 714.9722
$ws.Range("J74").Value = 2835.6667
$ws.Range("K74").Value = 714.9722
$ws.Range("L74").Value = 2835.6667
$ws.Range("M74").Value = 159.0278
$ws.Range("N74").Value = -4583.6667
$ws.Range("H77").Value = 1017.9286
$ws.Range("I77").Value = 714.9722
$ws.Range("J77").Value = 2835.6667
$ws.Range("K77").Value = 3574.861
$ws.Range("L77").Value = 14178.3335
$ws.Range("M77").Value = 793.1389999999997
$ws.Range("N77").Value = -22914.3335
$ws.Range("H94").Value = 24000
$ws.Range("J94").Value = 24000
$ws.Range("L94").Value = 24000
$ws.Range("N94").Value = -25802
$ws.Range("H110").Value = 1690.238
$ws.Range("I110").Value = 1290.8462
$ws.Range("J110").Value = 2339.25
$ws.Range("K110").Value = 1290.8462
$ws.Range("L110").Value = 2339.25
$ws.Range("M110").Value = 754.1538
$ws.Range("N110").Value = -6429.25
$ws.Range("H132").Value = 2490.8235
$ws.Range("I132").Value = 2781.4783
$ws.Range("J132").Value = 1883.091
$ws.Range("K132").Value = 8344.4349
$ws.Range("L132").Value = 5649.272999999999
$ws.Range("M132").Value = -5814.4349
$ws.Range("N132").Value = -10709.273
$ws.Range("H136").Value = 20000942
$ws.Range("I136").Value = 22223018
$ws.Range("J136").Value = 2262.8
$ws.Range("K136").Value = 66669054
$ws.Range("L136").Value = 6788.400000000001
$ws.Range("M136").Value = -66666504
$ws.Range("N136").Value = -11888.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1327
$ws.Range("H105").Value = 53153704
$ws.Range("I105").Value = 63119736
$ws.Range("K105").Value = 63119736
$ws.Range("M105").Value = -63117989
$ws.Range("H134").Value = 6601
$ws.Range("I134").Value = 2240.7
$ws.Range("J134").Value = 17501.75
$ws.Range("K134").Value = 6722.099999999999
$ws.Range("L134").Value = 52505.25
$ws.Range("M134").Value = -4187.099999999999
$ws.Range("N134").Value = -57575.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858830
$ws.Range("I16").Value = 166668270
$ws.Range("J16").Value = 2200
$ws.Range("K16").Value = 166668270
$ws.Range("L16").Value = 2200
$ws.Range("M16").Value = -166667983
$ws.Range("N16").Value = -2774
$ws.Range("H22").Value = 54342.383
$ws.Range("I22").Value = 718
$ws.Range("J22").Value = 87857.625
$ws.Range("K22").Value = 718
$ws.Range("L22").Value = 87857.625
$ws.Range("M22").Value = -368
$ws.Range("N22").Value = -88557.625
$ws.Range("H31").Value = 1761.8276
$ws.Range("I31").Value = 1853.0454
$ws.Range("K31").Value = 1853.0454
$ws.Range("M31").Value = -1558.0454
$ws.Range("H34").Value = 1761.8276
$ws.Range("I34").Value = 1853.0454
$ws.Range("K34").Value = 1853.0454
$ws.Range("M34").Value = -1651.0454
$ws.Range("H58").Value = 943.1923
$ws.Range("I58").Value = 843.3158
$ws.Range("J58").Value = 1214.2858
$ws.Range("K58").Value = 843.3158
$ws.Range("L58").Value = 1214.2858
$ws.Range("M58").Value = -640.3158
$ws.Range("N58").Value = -1620.2858
$ws.Range("H95").Value = 14350
$ws.Range("J95").Value = 14350
$ws.Range("L95").Value = 14350
$ws.Range("N95").Value = -19842
$ws.Range("H113").Value = 142858830
$ws.Range("I113").Value = 166668270
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 166668270
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -166666100
$ws.Range("N113").Value = -6540
$ws.Range("H132").Value = 2092.9473
$ws.Range("I132").Value = 1825.7059
$ws.Range("J132").Value = 4364.5
$ws.Range("K132").Value = 5477.1177
$ws.Range("L132").Value = 13093.5
$ws.Range("M132").Value = -2947.1177
$ws.Range("N132").Value = -18153.5
$ws.Range("H134").Value = 13158837
$ws.Range("I134").Value = 900.36
$ws.Range("J134").Value = 38462560
$ws.Range("K134").Value = 2701.08
$ws.Range("L134").Value = 115387680
$ws.Range("M134").Value = -166.0799999999999
$ws.Range("N134").Value = -115392750
$ws.Range("H136").Value = 943.1923
$ws.Range("I136").Value = 843.3158
$ws.Range("J136").Value = 1214.2858
$ws.Range("K136").Value = 2529.9474
$ws.Range("L136").Value = 3642.8574
$ws.Range("M136").Value = 20.05259999999998
$ws.Range("N136").Value = -8742.857400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2256.625
$ws.Range("I5").Value = 3017.3333
$ws.Range("J5").Value = 1278.5714
$ws.Range("K5").Value = 9051.999899999999
$ws.Range("L5").Value = 3835.7142
$ws.Range("M5").Value = -8939.999899999999
$ws.Range("N5").Value = -4059.7142
$ws.Range("H125").Value = 4539.3
$ws.Range("I125").Value = 2800
$ws.Range("J125").Value = 6278.6
$ws.Range("K125").Value = 8400
$ws.Range("L125").Value = 18835.8
$ws.Range("M125").Value = -3480
$ws.Range("N125").Value = -28675.8
$ws.Range("H126").Value = 5863.96
$ws.Range("I126").Value = 3399.8
$ws.Range("J126").Value = 6480
$ws.Range("K126").Value = 10199.4
$ws.Range("L126").Value = 19440
$ws.Range("M126").Value = -5259.400000000001
$ws.Range("N126").Value = -29320
$ws.Range("H131").Value = 18519784
$ws.Range("J131").Value = 1352.8776
$ws.Range("L131").Value = 4058.6328
$ws.Range("N131").Value = -14138.6328
$ws.Range("H133").Value = 3518.1924
$ws.Range("I133").Value = 1902.75
$ws.Range("J133").Value = 4236.1665
$ws.Range("K133").Value = 5708.25
$ws.Range("L133").Value = 12708.4995
$ws.Range("M133").Value = -648.25
$ws.Range("N133").Value = -22828.4995
$ws.Range("H135").Value = 2256.625
$ws.Range("I135").Value = 3017.3333
$ws.Range("J135").Value = 1278.5714
$ws.Range("K135").Value = 27155.9997
$ws.Range("L135").Value = 11507.1426
$ws.Range("M135").Value = -24620.9997
$ws.Range("N135").Value = -16577.1426

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2600.9333
$ws.Range("I80").Value = 1711
$ws.Range("J80").Value = 3045.9
$ws.Range("K80").Value = 1711
$ws.Range("L80").Value = 3045.9
$ws.Range("M80").Value = -713
$ws.Range("N80").Value = -5041.9
$ws.Range("H83").Value = 2600.9333
$ws.Range("I83").Value = 1711
$ws.Range("J83").Value = 3045.9
$ws.Range("K83").Value = 8555
$ws.Range("L83").Value = 15229.5
$ws.Range("M83").Value = -3563
$ws.Range("N83").Value = -25213.5
$ws.Range("H113").Value = 1285.8572
$ws.Range("I113").Value = 1602.75
$ws.Range("J113").Value = 863.3333
$ws.Range("K113").Value = 1602.75
$ws.Range("L113").Value = 863.3333
$ws.Range("M113").Value = 567.25
$ws.Range("N113").Value = -5203.3333
$ws.Range("H132").Value = 1567.0714
$ws.Range("I132").Value = 1435.12
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 4305.36
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -1775.36
$ws.Range("N132").Value = -13060.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 700005
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H100").Value = 763.3333
$ws.Range("I100").Value = 763.3333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1526.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -985.6666
$ws.Range("H113").Value = 351.1111
$ws.Range("I113").Value = 248.5
$ws.Range("J113").Value = 402.41666
$ws.Range("K113").Value = 745.5
$ws.Range("L113").Value = 1207.24998
$ws.Range("M113").Value = 1424.5
$ws.Range("N113").Value = -5547.249980000001
$ws.Range("H126").Value = 33334688
$ws.Range("I126").Value = 41667784
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 125003352
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -125000882
$ws.Range("N126").Value = -11840
$ws.Range("H132").Value = 3039.0645
$ws.Range("I132").Value = 3556.7083
$ws.Range("J132").Value = 1264.2858
$ws.Range("K132").Value = 10670.1249
$ws.Range("L132").Value = 3792.8574
$ws.Range("M132").Value = -8140.124899999999
$ws.Range("N132").Value = -8852.857400000001
$ws.Range("H136").Value = 621.7436
$ws.Range("J136").Value = 1152.1
$ws.Range("L136").Value = 3456.3
$ws.Range("N136").Value = -8556.299999999999
